# Adds 4 new city rows (Rejal Almaa, Shuwaq, Al Shaabah, Tumair) to the
# KSA Cities sheet, right after the existing last data row (row 135),
# matching the "Add files via upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 136; City = "Rejal Almaa"; Arabic = "رجال المع"; Lat = 17.407972000000001; Lon = 43.074370999999999; Area = "منطقة عسير";         Region = "جنوب المملكة" },
    @{ Row = 137; City = "Shuwaq";      Arabic = "شواق";       Lat = 19.919412999999999; Lon = 40.735523999999998; Area = "منطقة مكة المكرمة"; Region = "غرب المملكة" },
    @{ Row = 138; City = "Al Shaabah";  Arabic = "الشعبة";     Lat = 25.480855999999999; Lon = 49.622691000000003; Area = "المنطقة الشرقية";   Region = "شرق المملكة" },
    @{ Row = 139; City = "Tumair";      Arabic = "تمير";       Lat = 25.710602999999999; Lon = 45.872608;          Area = "منطقة الرياض";      Region = "وسط المملكة" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.City
    $ws.Cells.Item($row, 2).Value = $r.City
    $ws.Cells.Item($row, 3).Value = $r.Arabic
    $ws.Cells.Item($row, 4).Value = $r.Lat
    $ws.Cells.Item($row, 5).Value = $r.Lon
    $ws.Cells.Item($row, 6).Value = $r.Area
    $ws.Cells.Item($row, 7).Value = $r.Region
}

# Carry over the same cell formatting (thin borders etc.) used by every
# other data row, the same way copying row 135 down would in the UI.
$ws.Range("A135:G135").Copy()
$ws.Range("A136:G139").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Match the extended selection described in the commit's sheetView.
$ws.Range("A1:G139").Select()
